$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.55"
$ws.Range("E2").Value = "'0.89%"
$ws.Range("G2").Value = "'3"
$ws.Range("D3").Value = "'41.08"
$ws.Range("E3").Value = "'-0.75%"
$ws.Range("G3").Value = "'3"
$ws.Range("D4").Value = "'5.233"
$ws.Range("E4").Value = "'2.19%"
$ws.Range("G4").Value = "'3"
$ws.Range("D5").Value = "'0.07668"
$ws.Range("E5").Value = "'0.86%"
$ws.Range("G5").Value = "'3"
$ws.Range("D6").Value = "'1.625"
$ws.Range("E6").Value = "'-0.60%"
$ws.Range("G6").Value = "'3"
$ws.Range("D7").Value = "'0.9200"
$ws.Range("E7").Value = "'1.88%"
$ws.Range("G7").Value = "'3"
$ws.Range("D8").Value = "'2.474"
$ws.Range("E8").Value = "'-0.03%"
$ws.Range("G8").Value = "'3"
$ws.Range("D9").Value = "'0.1253"
$ws.Range("E9").Value = "'15.59%"
$ws.Range("G9").Value = "'3"
$ws.Range("D10").Value = "'0.1844"
$ws.Range("E10").Value = "'4.46%"
$ws.Range("G10").Value = "'3"
$ws.Range("D11").Value = "'0.09074"
$ws.Range("E11").Value = "'-1.37%"
$ws.Range("G11").Value = "'3"
$ws.Range("D12").Value = "'0.04358"
$ws.Range("E12").Value = "'1.64%"
$ws.Range("G12").Value = "'3"
$ws.Range("E13").Value = "'0.04%"
$ws.Range("G13").Value = "'3"
$ws.Range("E14").Value = "'1.05%"
$ws.Range("G14").Value = "'3"
$ws.Range("D15").Value = "'0.005792"
$ws.Range("E15").Value = "'-0.78%"
$ws.Range("G15").Value = "'3"
$ws.Range("E16").Value = "'2,392.75%"
$ws.Range("G16").Value = "'3"
$ws.Range("D17").Value = "'3.356"
$ws.Range("E17").Value = "'-0.15%"
$ws.Range("G17").Value = "'3"
$ws.Range("D18").Value = "'4.325"
$ws.Range("E18").Value = "'1.67%"
$ws.Range("G18").Value = "'3"
$ws.Range("G19").Value = "'3"
$ws.Range("D20").Value = "'7.216"
$ws.Range("E20").Value = "'9.36%"
$ws.Range("G20").Value = "'3"
$ws.Range("D21").Value = "'0.1383"
$ws.Range("E21").Value = "'1.37%"
$ws.Range("G21").Value = "'3"
$ws.Range("D22").Value = "'0.2925"
$ws.Range("E22").Value = "'9.06%"
$ws.Range("G22").Value = "'3"
$ws.Range("D23").Value = "'0.04058"
$ws.Range("E23").Value = "'-3.20%"
$ws.Range("G23").Value = "'3"
$ws.Range("D24").Value = "'0.001259"
$ws.Range("E24").Value = "'3.39%"
$ws.Range("G24").Value = "'3"
$ws.Range("D25").Value = "'0.004162"
$ws.Range("E25").Value = "'2.05%"
$ws.Range("G25").Value = "'3"
$ws.Range("D26").Value = "'0.0001272"
$ws.Range("E26").Value = "'-2.21%"
$ws.Range("G26").Value = "'3"
$ws.Range("G27").Value = "'3"
$ws.Range("G28").Value = "'3"
$ws.Range("G29").Value = "'3"
$ws.Range("G30").Value = "'3"
$ws.Range("G31").Value = "'3"
$ws.Range("G32").Value = "'3"
$ws.Range("G33").Value = "'3"
$ws.Range("G34").Value = "'3"
$ws.Range("G35").Value = "'3"
$ws.Range("G36").Value = "'3"
$ws.Range("G37").Value = "'3"
$ws.Range("D38").Value = "'0.02453"
$ws.Range("E38").Value = "'2.26%"
$ws.Range("G38").Value = "'3"
$ws.Range("D39").Value = "'0.05300"
$ws.Range("E39").Value = "'2.15%"
$ws.Range("G39").Value = "'3"
$ws.Range("D40").Value = "'0.007845"
$ws.Range("E40").Value = "'0.77%"
$ws.Range("G40").Value = "'3"
$ws.Range("D41").Value = "'0.1315"
$ws.Range("E41").Value = "'1.42%"
$ws.Range("G41").Value = "'3"
$ws.Range("G42").Value = "'3"
$ws.Range("E43").Value = "'-3.45%"
$ws.Range("G43").Value = "'3"
$ws.Range("D44").Value = "'0.008355"
$ws.Range("E44").Value = "'3.47%"
$ws.Range("G44").Value = "'3"
$ws.Range("D45").Value = "'0.3072"
$ws.Range("E45").Value = "'0.85%"
$ws.Range("G45").Value = "'3"
$ws.Range("D46").Value = "'0.00006667"
$ws.Range("E46").Value = "'-1.10%"
$ws.Range("G46").Value = "'3"
$ws.Range("E47").Value = "'0.13%"
$ws.Range("G47").Value = "'3"
$ws.Range("D48").Value = "'0.1699"
$ws.Range("E48").Value = "'1,532.28%"
$ws.Range("G48").Value = "'3"
$ws.Range("G49").Value = "'3"
$ws.Range("E50").Value = "'0.13%"
$ws.Range("G50").Value = "'3"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.13%"
$ws.Range("G51").Value = "'3"
